# Commiting Customer + Account files + ECRP
#
# Starting layout: Sheet1, Sheet0
# Target layout:   Sheet1, Sheet2 (new, empty), Sheet0, Sheet3 (new, ECRP data)
#
# Sheet1 also gets 4 new data rows (TC / Customer_ID / PD), and a new
# "Sheet3" tab carries the ECRP account rows.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet0 = $wb.Worksheets.Item("Sheet0")

# --- insert the two new (blank-for-now) worksheets in the right slots ---
$sheet2 = $wb.Worksheets.Add($null, $sheet1)   # lands right after Sheet1
$sheet3 = $wb.Worksheets.Add($null, $sheet0)   # lands right after Sheet0

# --- Sheet1: append 4 more customer rows (all text, like the existing row2) ---
function Set-TextCell($ws, $addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

$sheet1Rows = @(
    @("118463", "17705260", "1003"),
    @("118463", "17705261", "1007"),
    @("118464", "17705262", "1010"),
    @("118465", "17705263", "1011")
)

$rowIdx = 3
foreach ($row in $sheet1Rows) {
    Set-TextCell $sheet1 ("A" + $rowIdx) $row[0]
    Set-TextCell $sheet1 ("B" + $rowIdx) $row[1]
    Set-TextCell $sheet1 ("C" + $rowIdx) $row[2]
    $rowIdx++
}

[void]$sheet1.Range("A2").Select()
$sheet1.PageSetup.Orientation = 1

# --- Sheet3 (ECRP): header + 10 numeric data rows ---
$sheet3.Range("A1").Value = "TC"
$sheet3.Range("B1").Value = "Customer_ID"
$sheet3.Range("C1").Value = "PD"

$sheet3Rows = @(
    @(118463, 17705229, 1010),
    @(118463, 17705230, 6018),
    @(118463, 17705231, 1007),
    @(118464, 17705232, 1010),
    @(118465, 17705233, 1011),
    @(118466, 17705234, 6025),
    @(118468, 17705235, 6004),
    @(118469, 17705242, 6005),
    @(118469, 17705237, 1003),
    @(118471, 17705238, 1033)
)

$rowIdx = 2
foreach ($row in $sheet3Rows) {
    $sheet3.Range("A" + $rowIdx).Value = $row[0]
    $sheet3.Range("B" + $rowIdx).Value = $row[1]
    $sheet3.Range("C" + $rowIdx).Value = $row[2]
    $rowIdx++
}

[void]$sheet3.Range("A2:C11").Select()

# Keep Sheet1 as the active/visible tab, matching the original workbook.
$sheet1.Activate()
